$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# New long text blocks for the "Sprint 2" table (rows 19-34)
# ---------------------------------------------------------------
$sSiteNuvem = @"
Site Institucional na Nuvem
Cadastro, Login e Dashboard na Nuvem,                               conectado com BD                                                                                                     Mapeamento das tabelas (entidades) em classes Javascript
"@
$sFluxograma = @"
Fluxograma do Processo de Atendimento do Suporte
Ferramenta de Help Desk configurada e integrada à solução
"@
$sModelagem = @"
                                                                                                          Modelagem Lógica (Final)
Script de criação do Banco (Final)
Tabelas criadas no Azure (Final)

"@
$sTesteIntegrado = @"
Teste Integrado do Analytics
Teste Integrado da Solução de IoT                                         (Simulador + Banco de Dados)
"@
$sManual = @"
Manual de Instalação
 Doc. Final do Projeto
 PPT da Apresentação do Projeto
 Prévia (Demonstração da Solução + Apresentação)
"@
$sEncarregado = @"
Será feito, estudado e 
realizado por todos.
"@

# ---------------------------------------------------------------
# 1) Copy the formatting (borders/fills/number-format/alignment)
#    from the existing Sprint-1 table down onto the new block so
#    every new cell inherits the same per-column styling.
# ---------------------------------------------------------------
$ws.Range("A1:D16").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 2) Header row (row 19) - same headers as row 1
# ---------------------------------------------------------------
$ws.Range("A19").Value = "Matérias"
$ws.Range("B19").Value = "Entregáveis"
$ws.Range("C19").Value = "Encarregado"
$ws.Range("D19").Value = "Data da Sprint"

# ---------------------------------------------------------------
# 3) Data rows
# ---------------------------------------------------------------
$ws.Range("A20").Value = "Algoritmos"
$ws.Range("B20").Value = $sSiteNuvem
$ws.Range("C20").Value = $sEncarregado
$ws.Range("D20").Value = 44357

$ws.Range("A24").Value = "Arq. Comp"
$ws.Range("B24").Value = $sTesteIntegrado
$ws.Range("C24").Value = $sEncarregado
$ws.Range("D24").Value = 44357

$ws.Range("A27").Value = "Banco de dados"
$ws.Range("B27").Value = $sModelagem
$ws.Range("C27").Value = $sEncarregado
$ws.Range("D27").Value = 44357

$ws.Range("A29").Value = "PI"
$ws.Range("B29").Value = $sManual
$ws.Range("C29").Value = $sEncarregado
$ws.Range("D29").Value = 44357

$ws.Range("A32").Value = "TI"
$ws.Range("B32").Value = $sFluxograma
$ws.Range("C32").Value = $sEncarregado
$ws.Range("D32").Value = 44357

# ---------------------------------------------------------------
# 4) Merges (mirrors the 5 groups used by the Sprint-1 table)
# ---------------------------------------------------------------
$ws.Range("A20:A23").Merge()
$ws.Range("B20:B23").Merge()
$ws.Range("C20:C23").Merge()
$ws.Range("D20:D23").Merge()

$ws.Range("A24:A26").Merge()
$ws.Range("B24:B26").Merge()
$ws.Range("C24:C26").Merge()
$ws.Range("D24:D26").Merge()

$ws.Range("A27:A28").Merge()
$ws.Range("B27:B28").Merge()
$ws.Range("C27:C28").Merge()
$ws.Range("D27:D28").Merge()

$ws.Range("A29:A31").Merge()
$ws.Range("B29:B31").Merge()
$ws.Range("C29:C31").Merge()
$ws.Range("D29:D31").Merge()

$ws.Range("A32:A34").Merge()
$ws.Range("B32:B34").Merge()
$ws.Range("C32:C34").Merge()
$ws.Range("D32:D34").Merge()

# ---------------------------------------------------------------
# 5) Row heights that Excel auto-fit once the long, wrapped text
#    was entered into the merged cells.
# ---------------------------------------------------------------
$ws.Rows.Item(23).RowHeight = 29.25
$ws.Rows.Item(27).RowHeight = 39
$ws.Rows.Item(28).RowHeight = 18.75
$ws.Rows.Item(31).RowHeight = 42

# ---------------------------------------------------------------
# 6) Leave the selection/scroll where the author left it
# ---------------------------------------------------------------
$ws.Range("J33").Select()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
